$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Lateral Movement" row (row 15): Local went Fully -> Partially,
# Network went No -> Partially (multi-facility effect removal now partially
# works across the network).
$ws.Range("B15").Value = "Partially"
$ws.Range("C15").Value = "Partially"

# Update the "Remove all backdoor effects" row (row 19): Local went No -> Partially.
$ws.Range("B19").Value = "Partially"

# Move the active selection on the sheet to D32 (was C33).
$ws.Range("D32").Select()
